{"js": "// Update the date and all multiplication problems to the new \"output\n// generated at c986bee\" values.\nconst replacements = [\n  [\"2024-12-30 Monday\", \"2024-12-31 Tuesday\"],\n  [\"297\u00d77=\", \"406\u00d78=\"],\n  [\"402\u00d78=\", \"853\u00d73=\"],\n  [\"503\u00d75=\", \"147\u00d79=\"],\n  [\"294\u00d77=\", \"797\u00d76=\"],\n  [\"883\u00d78=\", \"322\u00d72=\"],\n  [\"160\u00d72=\", \"359\u00d74=\"],\n  [\"511\u00d76=\", \"408\u00d76=\"],\n  [\"656\u00d78=\", \"761\u00d74=\"],\n  [\"335\u00d77=\", \"365\u00d79=\"],\n  [\"498\u00d72=\", \"565\u00d75=\"],\n  [\"878\u00d76=\", \"951\u00d72=\"],\n  [\"272\u00d75=\", \"613\u00d79=\"],\n  [\"145\u00d78=\", \"283\u00d76=\"],\n  [\"743\u00d78=\", \"537\u00d78=\"],\n  [\"815\u00d74=\", \"747\u00d75=\"],\n  [\"954\u00d72=\", \"268\u00d76=\"],\n  [\"155\u00d79=\", \"944\u00d76=\"],\n  [\"729\u00d76=\", \"468\u00d79=\"],\n  [\"605\u00d79=\", \"304\u00d77=\"],\n  [\"467\u00d77=\", \"465\u00d79=\"],\n  [\"497\u00d78=\", \"531\u00d76=\"],\n  [\"722\u00d78=\", \"966\u00d74=\"],\n  [\"343\u00d72=\", \"822\u00d75=\"],\n  [\"227\u00d74=\", \"213\u00d74=\"],\n  [\"455\u00d75=\", \"607\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and all multiplication problems to the new \"output\n# generated at c986bee\" values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-30 Monday\", \"2024-12-31 Tuesday\"),\n    @(\"297\u00d77=\", \"406\u00d78=\"),\n    @(\"402\u00d78=\", \"853\u00d73=\"),\n    @(\"503\u00d75=\", \"147\u00d79=\"),\n    @(\"294\u00d77=\", \"797\u00d76=\"),\n    @(\"883\u00d78=\", \"322\u00d72=\"),\n    @(\"160\u00d72=\", \"359\u00d74=\"),\n    @(\"511\u00d76=\", \"408\u00d76=\"),\n    @(\"656\u00d78=\", \"761\u00d74=\"),\n    @(\"335\u00d77=\", \"365\u00d79=\"),\n    @(\"498\u00d72=\", \"565\u00d75=\"),\n    @(\"878\u00d76=\", \"951\u00d72=\"),\n    @(\"272\u00d75=\", \"613\u00d79=\"),\n    @(\"145\u00d78=\", \"283\u00d76=\"),\n    @(\"743\u00d78=\", \"537\u00d78=\"),\n    @(\"815\u00d74=\", \"747\u00d75=\"),\n    @(\"954\u00d72=\", \"268\u00d76=\"),\n    @(\"155\u00d79=\", \"944\u00d76=\"),\n    @(\"729\u00d76=\", \"468\u00d79=\"),\n    @(\"605\u00d79=\", \"304\u00d77=\"),\n    @(\"467\u00d77=\", \"465\u00d79=\"),\n    @(\"497\u00d78=\", \"531\u00d76=\"),\n    @(\"722\u00d78=\", \"966\u00d74=\"),\n    @(\"343\u00d72=\", \"822\u00d75=\"),\n    @(\"227\u00d74=\", \"213\u00d74=\"),\n    @(\"455\u00d75=\", \"607\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # 1 = wdFindContinue, 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
